$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_19")

# Re-arrange the footnote block:
#  - B11 now carries the "Actualización" note (reworded, trailing period)
#  - B12 now carries the "Fuente" note (trailing period added)
#  - F12 / F13 are cleared out (they used to hold these two notes)
$ws.Range("B11").Value = "Actualización: mayo 2024."
$ws.Range("B12").Value = "Fuente: SICT. Subsecretaria de Comunicaciones y Transportes. Dirección General de Autotransporte Federal."
$ws.Range("F12").Value = ""
$ws.Range("F13").Value = ""

# Fix accented "Kilómetros" in the two header cells (D4, F4)
$ws.Range("D4").Value = "Pasajeros-Kilómetros transportados (Millones)"
$ws.Range("F4").Value = "Toneladas-kilómetros transportadas (Millones)"
